# The cancelled event "南宁·0316全职only-全明星周末（取消）" (row 2) was removed
# from both the "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet.
# Deleting the row shifts every following row up by one; a handful of the
# remaining events also picked up slightly higher "想去人数" (interest count)
# numbers in column F in this later scrape snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Rows.Item(2).Delete()

# Renumber the sequence column (A) now that row 2 is gone.
$wsExpo.Range("A2").Value = 1
$wsExpo.Range("A3").Value = 2
$wsExpo.Range("A4").Value = 3
$wsExpo.Range("A5").Value = 4

# Updated "想去人数" counts for the events that remain.
$wsExpo.Range("F2").Value = 2184
$wsExpo.Range("F3").Value = 896
$wsExpo.Range("F4").Value = 1565
$wsExpo.Range("F5").Value = 379

# ---- Sheet "全部类型" (All types) ----
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Rows.Item(2).Delete()

# Renumber the sequence column (A) now that row 2 is gone.
$wsAll.Range("A2").Value = 1
$wsAll.Range("A3").Value = 2
$wsAll.Range("A4").Value = 3
$wsAll.Range("A5").Value = 4
$wsAll.Range("A6").Value = 5
$wsAll.Range("A7").Value = 6

# Updated "想去人数" counts for the events that remain.
$wsAll.Range("F2").Value = 2184
$wsAll.Range("F5").Value = 896
$wsAll.Range("F6").Value = 1565
$wsAll.Range("F7").Value = 379
